$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2291666666666667
$ws.Range("C2").Value = 0.5
$ws.Range("J2").Value = 0.02083333333333333
$ws.Range("P2").Value = 0.1458333333333333
$ws.Range("S2").Value = 0.1041666666666667

# Row 3
$ws.Range("J3").Value = 0.04
$ws.Range("P3").Value = 0.64
$ws.Range("S3").Value = 0.32

# Row 4
$ws.Range("O4").Value = 0.1111111111111111
$ws.Range("P4").Value = 0.5555555555555556

# Row 6
$ws.Range("B6").Value = 0.03448275862068965
$ws.Range("J6").Value = 0.2758620689655172
$ws.Range("O6").Value = 0.1379310344827586
$ws.Range("Q6").Value = 0.103448275862069
$ws.Range("R6").Value = 0.103448275862069
$ws.Range("S6").Value = 0.3448275862068966

# Row 7
$ws.Range("B7").Value = 0.1
$ws.Range("F7").Value = 0.2
$ws.Range("J7").Value = 0.1
$ws.Range("Q7").Value = 0.2
$ws.Range("R7").Value = 0.1
$ws.Range("S7").Value = 0.3

# Row 8
$ws.Range("B8").Value = 0.225
$ws.Range("F8").Value = 0.15
$ws.Range("J8").Value = 0.05
$ws.Range("Q8").Value = 0.125
$ws.Range("R8").Value = 0.025
$ws.Range("S8").Value = 0.425

# Row 9
$ws.Range("B9").Value = 0.2307692307692308
$ws.Range("D9").Value = 0.03846153846153846
$ws.Range("F9").Value = 0.03846153846153846
$ws.Range("J9").Value = 0.07692307692307693
$ws.Range("Q9").Value = 0.4615384615384616
$ws.Range("R9").Value = 0.03846153846153846
$ws.Range("S9").Value = 0.1153846153846154

# Row 10
$ws.Range("B10").Value = 0.1104972375690608
$ws.Range("D10").Value = 0.03867403314917127
$ws.Range("F10").Value = 0.05524861878453038
$ws.Range("J10").Value = 0.1270718232044199
$ws.Range("O10").Value = 0.02762430939226519
$ws.Range("Q10").Value = 0.2099447513812155
$ws.Range("R10").Value = 0.04972375690607735
$ws.Range("S10").Value = 0.3812154696132597

# Row 11
$ws.Range("G11").Value = 0.1304347826086956
$ws.Range("J11").Value = 0.1739130434782609
$ws.Range("K11").Value = 0.1739130434782609
$ws.Range("L11").Value = 0.5217391304347826

# Row 12
$ws.Range("G12").Value = 0.5
$ws.Range("J12").Value = 0.5

# Row 13
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.5

# Row 15
$ws.Range("F15").Value = 0.06451612903225806
$ws.Range("H15").Value = 0.09677419354838709
$ws.Range("I15").Value = 0.06451612903225806
$ws.Range("J15").Value = 0.3870967741935484
$ws.Range("K15").Value = 0.03225806451612903
$ws.Range("M15").Value = 0.03225806451612903
$ws.Range("S15").Value = 0.3225806451612903

# Row 16
$ws.Range("F16").Value = 0.03703703703703703
$ws.Range("H16").Value = 0.1851851851851852
$ws.Range("I16").Value = 0.07407407407407407
$ws.Range("J16").Value = 0.6296296296296297
$ws.Range("K16").Value = 0.03703703703703703
$ws.Range("S16").Value = 0.03703703703703703

# Row 17
$ws.Range("H17").Value = 0.06779661016949153
$ws.Range("I17").Value = 0.1355932203389831
$ws.Range("J17").Value = 0.6101694915254238
$ws.Range("K17").Value = 0.0847457627118644
$ws.Range("O17").Value = 0.05084745762711865
$ws.Range("S17").Value = 0.05084745762711865

# Row 18
$ws.Range("F18").Value = 0.06666666666666667
$ws.Range("H18").Value = 0.1333333333333333
$ws.Range("I18").Value = 0.06666666666666667
$ws.Range("J18").Value = 0.5333333333333333
$ws.Range("O18").Value = 0.1333333333333333
$ws.Range("S18").Value = 0.06666666666666667

# Row 19
$ws.Range("F19").Value = 0.006944444444444444
$ws.Range("H19").Value = 0.1736111111111111
$ws.Range("I19").Value = 0.09027777777777778
$ws.Range("J19").Value = 0.4513888888888889
$ws.Range("K19").Value = 0.0763888888888889
$ws.Range("M19").Value = 0.006944444444444444
$ws.Range("O19").Value = 0.1041666666666667
$ws.Range("S19").Value = 0.09027777777777778
